# chore: update Sheets via scheduled runner
# Refreshes cached market-price/profit columns (H:N) on a handful of leve
# rows across each job sheet, pulling in newer currentAveragePrice* figures
# and recomputed LevePriceNQ/HQ + profit deltas.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3920.75
$ws.Range("I43").Value = 4660.3335
$ws.Range("J43").Value = 1702
$ws.Range("K43").Value = 4660.3335
$ws.Range("L43").Value = 1702
$ws.Range("M43").Value = -4591.3335
$ws.Range("N43").Value = -1840
$ws.Range("H53").Value = 462.15
$ws.Range("I53").Value = 439.9375
$ws.Range("K53").Value = 439.9375
$ws.Range("M53").Value = 197.0625
$ws.Range("H58").Value = 1509972.9
$ws.Range("I58").Value = 3268358
$ws.Range("J58").Value = 2785.7144
$ws.Range("K58").Value = 9805074
$ws.Range("L58").Value = 8357.143199999999
$ws.Range("M58").Value = -9804924
$ws.Range("N58").Value = -8657.143199999999
$ws.Range("H64").Value = 79800
$ws.Range("I64").Value = 252250
$ws.Range("J64").Value = 3155.5557
$ws.Range("K64").Value = 252250
$ws.Range("L64").Value = 3155.5557
$ws.Range("M64").Value = -252002
$ws.Range("N64").Value = -3651.5557
$ws.Range("H67").Value = 79800
$ws.Range("I67").Value = 252250
$ws.Range("J67").Value = 3155.5557
$ws.Range("K67").Value = 252250
$ws.Range("L67").Value = 3155.5557
$ws.Range("M67").Value = -251392
$ws.Range("N67").Value = -4871.5557
$ws.Range("I113").Value = 252001.25
$ws.Range("K113").Value = 252001.25
$ws.Range("M113").Value = -248747.25
$ws.Range("H116").Value = 33333
$ws.Range("I116").Value = 33333
$ws.Range("K116").Value = 33333
$ws.Range("M116").Value = -29891
$ws.Range("H129").Value = 2499.7424
$ws.Range("I129").Value = 8246.23
$ws.Range("J129").Value = 1090.2264
$ws.Range("K129").Value = 24738.69
$ws.Range("L129").Value = 3270.6792
$ws.Range("M129").Value = -19738.69
$ws.Range("N129").Value = -13270.6792
$ws.Range("H138").Value = 2797.9343
$ws.Range("I138").Value = 1631.6538
$ws.Range("J138").Value = 3664.3142
$ws.Range("K138").Value = 4894.9614
$ws.Range("L138").Value = 10992.9426
$ws.Range("M138").Value = 245.0385999999999
$ws.Range("N138").Value = -21272.9426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26531.402
$ws.Range("I32").Value = 7054.028
$ws.Range("K32").Value = 7054.028
$ws.Range("M32").Value = -6767.028
$ws.Range("H74").Value = 962.72
$ws.Range("I74").Value = 888.4211
$ws.Range("J74").Value = 1198
$ws.Range("K74").Value = 888.4211
$ws.Range("L74").Value = 1198
$ws.Range("M74").Value = -14.42110000000002
$ws.Range("N74").Value = -2946
$ws.Range("H77").Value = 962.72
$ws.Range("I77").Value = 888.4211
$ws.Range("J77").Value = 1198
$ws.Range("K77").Value = 4442.1055
$ws.Range("L77").Value = 5990
$ws.Range("M77").Value = -74.10549999999967
$ws.Range("N77").Value = -14726
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H97").Value = 40117.348
$ws.Range("I97").Value = 44689.086
$ws.Range("J97").Value = 5067.3335
$ws.Range("K97").Value = 44689.086
$ws.Range("L97").Value = 5067.3335
$ws.Range("M97").Value = -44193.086
$ws.Range("N97").Value = -6059.3335
$ws.Range("H110").Value = 143157980
$ws.Range("I110").Value = 143157980
$ws.Range("K110").Value = 143157980
$ws.Range("M110").Value = -143155935
$ws.Range("H132").Value = 13535.2
$ws.Range("I132").Value = 15826.3
$ws.Range("K132").Value = 47478.89999999999
$ws.Range("M132").Value = -44948.89999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 67016.82000000001
$ws.Range("I86").Value = 86715.84
$ws.Range("J86").Value = 2995
$ws.Range("K86").Value = 86715.84
$ws.Range("L86").Value = 2995
$ws.Range("M86").Value = -85592.84
$ws.Range("N86").Value = -5241
$ws.Range("H89").Value = 67016.82000000001
$ws.Range("I89").Value = 86715.84
$ws.Range("J89").Value = 2995
$ws.Range("K89").Value = 433579.2
$ws.Range("L89").Value = 14975
$ws.Range("M89").Value = -427963.2
$ws.Range("N89").Value = -26207
$ws.Range("H94").Value = 753.5238000000001
$ws.Range("I94").Value = 740.3889
$ws.Range("J94").Value = 832.3333
$ws.Range("K94").Value = 740.3889
$ws.Range("L94").Value = 832.3333
$ws.Range("M94").Value = -289.3889
$ws.Range("N94").Value = -1734.3333
$ws.Range("H99").Value = 1942.2222
$ws.Range("I99").Value = 1726.6666
$ws.Range("J99").Value = 1985.3334
$ws.Range("K99").Value = 1726.6666
$ws.Range("L99").Value = 1985.3334
$ws.Range("M99").Value = -228.6666
$ws.Range("N99").Value = -4981.3334
$ws.Range("H134").Value = 12419.92
$ws.Range("I134").Value = 12978.195
$ws.Range("K134").Value = 38934.585
$ws.Range("M134").Value = -36399.585

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 189.86667
$ws.Range("I7").Value = 111.5
$ws.Range("J7").Value = 279.42856
$ws.Range("K7").Value = 111.5
$ws.Range("L7").Value = 279.42856
$ws.Range("M7").Value = 1.5
$ws.Range("N7").Value = -505.42856
$ws.Range("H134").Value = 1873.8667
$ws.Range("I134").Value = 835
$ws.Range("K134").Value = 2505
$ws.Range("M134").Value = 30

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 407.3684
$ws.Range("I23").Value = 17
$ws.Range("J23").Value = 480.5625
$ws.Range("K23").Value = 51
$ws.Range("L23").Value = 1441.6875
$ws.Range("M23").Value = 184
$ws.Range("N23").Value = -1911.6875
$ws.Range("H33").Value = 2567.8333
$ws.Range("I33").Value = 2038.8
$ws.Range("J33").Value = 2945.7144
$ws.Range("K33").Value = 12232.8
$ws.Range("L33").Value = 17674.2864
$ws.Range("M33").Value = -11949.8
$ws.Range("N33").Value = -18240.2864
$ws.Range("H34").Value = 1217.7778
$ws.Range("J34").Value = 1288.2354
$ws.Range("L34").Value = 3864.7062
$ws.Range("N34").Value = -4032.7062

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 36074
$ws.Range("I127").Value = 30148
$ws.Range("J127").Value = 42000
$ws.Range("K127").Value = 30148
$ws.Range("L127").Value = 42000
$ws.Range("M127").Value = -25188
$ws.Range("N127").Value = -51920
$ws.Range("H132").Value = 2849.484
$ws.Range("I132").Value = 2068.65
$ws.Range("J132").Value = 4269.1816
$ws.Range("K132").Value = 6205.950000000001
$ws.Range("L132").Value = 12807.5448
$ws.Range("M132").Value = -3675.950000000001
$ws.Range("N132").Value = -17867.5448

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 30150
$ws.Range("J6").Value = 30150
$ws.Range("L6").Value = 30150
$ws.Range("N6").Value = -30374
$ws.Range("H22").Value = 1527.3334
$ws.Range("J22").Value = 559.3333
$ws.Range("L22").Value = 559.3333
$ws.Range("N22").Value = -1149.3333
$ws.Range("H27").Value = 1527.3334
$ws.Range("J27").Value = 559.3333
$ws.Range("L27").Value = 559.3333
$ws.Range("N27").Value = -773.3333
$ws.Range("H93").Value = 1974.129
$ws.Range("I93").Value = 1963.5238
$ws.Range("J93").Value = 1996.4
$ws.Range("K93").Value = 1963.5238
$ws.Range("L93").Value = 1996.4
$ws.Range("M93").Value = -715.5237999999999
$ws.Range("N93").Value = -4492.4
$ws.Range("H132").Value = 3060.7932
$ws.Range("I132").Value = 3290.8845
$ws.Range("K132").Value = 9872.6535
$ws.Range("M132").Value = -7342.6535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H96").Value = 111112730
$ws.Range("I96").Value = 166668660
$ws.Range("J96").Value = 866.6667
$ws.Range("K96").Value = 166668660
$ws.Range("L96").Value = 866.6667
$ws.Range("M96").Value = -166667287
$ws.Range("N96").Value = -3612.6667
$ws.Range("H136").Value = 17124.986
$ws.Range("I136").Value = 33733.8
$ws.Range("K136").Value = 101201.4
$ws.Range("M136").Value = -98651.40000000001
